$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Committee members" text: add a missing space before "Camilo".
$ws.Range("E3").Value = "Committee members: Sunshine Van Bael, Ph.D. (adviser; dissertation chair), Kathleen Ferris, Ph.D.  (co-advisor), Keith Clay, Ph.D., & P. Camilo Zalamea, Ph.D."

# Update the active selection to E3.
$null = $ws.Range("E3").Select()
